# Append a new task row (row 3) to the "Task Data" sheet, mirroring the
# layout/formatting already used by the existing data row (row 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 3

$ws.Cells.Item($newRow, 1).Value = "a"
$ws.Cells.Item($newRow, 2).Value = "b"
$ws.Cells.Item($newRow, 3).Value = 45208.48451388889
$ws.Cells.Item($newRow, 4).Value = 45208.48454861111
$ws.Cells.Item($newRow, 5).Value = 294.004

# Match the Start/End DateTime formatting used by the previous row so the
# new cells pick up the same style (rather than the default/general one).
$dateFormat = $ws.Cells.Item(2, 3).NumberFormat
$ws.Cells.Item($newRow, 3).NumberFormat = $dateFormat
$ws.Cells.Item($newRow, 4).NumberFormat = $dateFormat
